$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("huc12")

# Insert a new column before column D to hold the HUC12 drainage area (sq km).
# This shifts all the existing land-cover columns (old D..R) right by one
# (new E..S) and Excel auto-adjusts the formulas to match.
$ws.Columns("D").Insert()

# Header for the new column
$ws.Range("D1").Value = "huc12_area_sq.km"

# huc12 area values (sq km) for each site (rows without a known area are left blank)
$ws.Range("D3").Value = 215.36
$ws.Range("D4").Value = 107.5
$ws.Range("D5").Value = 200.22
$ws.Range("D6").Value = 130.85
$ws.Range("D7").Value = 90.75
$ws.Range("D9").Value = 86.9
$ws.Range("D10").Value = 133.66

# Match the column formatting used by the rest of the wide text columns
$ws.Columns("D").ColumnWidth = 46.5

# Restore/update the active selection like the saved workbook shows
$ws.Range("D10").Select()
